$wb = $excel.ActiveWorkbook

# The new sheet lives where "Sheet3" used to be (rId3 / sheetId 3).
$ws3 = $wb.Worksheets.Item("Sheet3")

# --- Column A: a few stand-alone formulas -------------------------------
$ws3.Range("A1").Formula = "=SUM(Unformatted!B2:B8)"
$ws3.Range("A2").Formula = "=A1*2"
$ws3.Range("A3").Formula = "=SUM(A1:A2)"

# --- Column D: literal values 1..13 -------------------------------------
for ($i = 1; $i -le 13; $i++) {
    $ws3.Cells.Item($i, 4).Value = $i
}

# --- Column E: D*10 (single formula row 1, shared formula rows 2-13) ----
$ws3.Range("E1").Formula = "=D1*10"
$ws3.Range("E2:E13").Formula = "=D2*10"

# --- Column F: SUM(D:E) (single formula row 1, shared formula rows 2-13)-
$ws3.Range("F1").Formula = "=SUM(D1:E1)"
$ws3.Range("F2:F13").Formula = "=SUM(D2:E2)"

# --- G13: grand total -----------------------------------------------------
$ws3.Range("G13").Formula = "=SUM(E13:F13)"

# Rename the sheet to reflect its new purpose.
$ws3.Name = "Formula"

# Make it the active/selected sheet (moves tabSelected + activeTab here,
# and clears tabSelected on the previously-active "Simple Format" sheet).
$ws3.Activate()
$null = $ws3.Range("G14").Select()
